$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct country name mappings
$ws.Range("A21").Value = "The Gambia"
$ws.Range("A25").Value = "Côte d'Ivoire"
$ws.Range("A40").Value = "Republic of Congo"

# Clear gdp / gdp_per_capita values that are no longer available,
# keeping the cells present (empty) rather than removing them entirely.
$cells = @("E21", "F21", "E25", "F25", "E40", "F40")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.ClearContents()
    $rng.Style = "Normal"
}
